$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

$ws.Range('D2').Value = '75.597.53'
$ws.Range('E2').Value = '  +8.90%  '
$ws.Range('D3').Value = '2.720.34'
$ws.Range('E3').Value = '  +12.21%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '187.14'
$ws.Range('E5').Value = '  +12.24%  '
Set-TextValue 'D6' '591.48'
$ws.Range('E6').Value = '  +4.95%  '
$ws.Range('E7').Value = '  -0.13%  '
Set-TextValue 'D8' '0.542'
$ws.Range('E8').Value = '  +5.53%  '
Set-TextValue 'D9' '0.196'
$ws.Range('E9').Value = '  +15.25%  '
$ws.Range('D10').Value = '2.719.17'
$ws.Range('E10').Value = '  +12.06%  '
$ws.Range('E11').Value = '  +1.36%  '
Set-TextValue 'D12' '0.363'
$ws.Range('E12').Value = '  +8.79%  '
Set-TextValue 'D13' '4.79'
$ws.Range('E13').Value = '  +2.60%  '
$ws.Range('D14').Value = '3.218.38'
$ws.Range('E14').Value = '  +12.11%  '
$ws.Range('D15').Value = '75.411.28'
$ws.Range('E15').Value = '  +8.77%  '
$ws.Range('E16').Value = '  +6.63%  '
Set-TextValue 'D17' '27.01'
$ws.Range('E17').Value = '  +12.58%  '
$ws.Range('D18').Value = '2.707.77'
$ws.Range('E18').Value = '  +11.80%  '
Set-TextValue 'D19' '9.33'
$ws.Range('E19').Value = '  +29.81%  '
Set-TextValue 'D20' '12.16'
$ws.Range('E20').Value = '  +12.11%  '
Set-TextValue 'D21' '378.13'
$ws.Range('E21').Value = '  +10.24%  '
Set-TextValue 'D22' '2.30'
$ws.Range('E22').Value = '  +14.77%  '
Set-TextValue 'D23' '4.12'
$ws.Range('E23').Value = '  +6.76%  '
$ws.Range('E24').Value = '  +4.85%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D25' '0.999'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D26' '71.00'
$ws.Range('E26').Value = '  +8.14%  '
Set-TextValue 'D27' '4.23'
$ws.Range('E27').Value = '  +11.03%  '
Set-TextValue 'D28' '9.59'
$ws.Range('E28').Value = '  +13.22%  '
$ws.Range('D29').Value = '2.858.01'
$ws.Range('E29').Value = '  +11.91%  '
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('D31').Value = '0.0₃0993'
$ws.Range('E31').Value = '  +16.09%  '
Set-TextValue 'D32' '523.78'
$ws.Range('E32').Value = '  +14.41%  '
Set-TextValue 'D33' '1.42'
$ws.Range('E33').Value = '  +13.01%  '
Set-TextValue 'D34' '7.88'
$ws.Range('E34').Value = '  +7.05%  '
$ws.Range('E35').Value = '  +10.91%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +8.08%  '
Set-TextValue 'D38' '161.17'
$ws.Range('E38').Value = '  +1.51%  '
Set-TextValue 'D39' '19.69'
$ws.Range('E39').Value = '  +8.03%  '
$ws.Range('E40').Value = '  +1.43%  '
$ws.Range('E41').Value = '  -0.02%  '
Set-TextValue 'D42' '173.56'
$ws.Range('E42').Value = '  +27.93%  '
$ws.Range('E43').Value = '  +14.53%  '
$ws.Range('E44').Value = '  +13.13%  '
$ws.Range('E45').Value = '  +10.14%  '
Set-TextValue 'D46' '1.23'
$ws.Range('E46').Value = '  +13.70%  '
Set-TextValue 'D47' '2.41'
$ws.Range('E47').Value = '  +14.62%  '
$ws.Range('E48').Value = '  +3.12%  '
$ws.Range('E49').Value = '  +18.59%  '
Set-TextValue 'D50' '3.72'
$ws.Range('E50').Value = '  +9.64%  '
$ws.Range('E51').Value = '  +11.88%  '
